# Update gh-pages output (generated data refresh) across all sheets.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Plain "want to go" count (column F) refreshes for rows unaffected by the
# later row insertion (rows 2-23).
$ws1.Range("F2").Value = 1272
$ws1.Range("F5").Value = 6479
$ws1.Range("F6").Value = 1812
$ws1.Range("F7").Value = 6382
$ws1.Range("F8").Value = 144
$ws1.Range("F9").Value = 1947
$ws1.Range("F11").Value = 18
$ws1.Range("F13").Value = 37
$ws1.Range("F17").Value = 8004
$ws1.Range("F20").Value = 188
$ws1.Range("F22").Value = 1753
$ws1.Range("F23").Value = 849

# Insert a brand-new row 29 ("杭州·次元幻想..."), pushing the former rows
# 29-37 down to 30-38.
$ws1.Rows.Item(29).Insert()

# New row 29 needs the same look (border / bold / centered) as the rest of
# column A's index cells.
$idxCell = $ws1.Range("A29")
$idxCell.Borders.LineStyle = 1
$idxCell.HorizontalAlignment = -4108
$idxCell.VerticalAlignment = -4160
$idxCell.Font.Bold = $true
$idxCell.Value = 28

# B29 holds a date-shaped label that must stay plain text (not get
# auto-converted to a date serial by Excel's input parsing).
$ws1.Range("B29").NumberFormat = "@"
$ws1.Range("B29").Value = "2024-10-05"

$ws1.Range("C29").Value = "杭州·次元幻想【玩美大舞台&吃/换谷大会】（免费活动）"
$ws1.Range("D29").Value = "文三路 玩美的一天沉浸式生活街区"
$ws1.Range("E29").Value = "2024.10.05 10:00-10.05 17:00"
$ws1.Range("F29").Value = 0
$ws1.Range("G29").Value = 30
$ws1.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=92028"
$ws1.Range("I29").Value = "//i2.hdslb.com/bfs/openplatform/202409/FaEB96xH1725394323651.jpeg"

# The "No." column (A) is a plain sequential index, not a formula, so the
# rows that shifted down keep their old literal number; bump them back in
# line with their new row position.
$ws1.Range("A30").Value = 29
$ws1.Range("A31").Value = 30
$ws1.Range("A32").Value = 31
$ws1.Range("A33").Value = 32
$ws1.Range("A34").Value = 33
$ws1.Range("A35").Value = 34
$ws1.Range("A36").Value = 35
$ws1.Range("A37").Value = 36
$ws1.Range("A38").Value = 37

# A few of the rows that shifted down also picked up refreshed counts.
$ws1.Range("F30").Value = 1812
$ws1.Range("F31").Value = 815
$ws1.Range("F32").Value = 387
$ws1.Range("F35").Value = 16

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 370
$ws2.Range("F5").Value = 206
$ws2.Range("F18").Value = 91
$ws2.Range("F22").Value = 34

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9546
$ws3.Range("F3").Value = 2287
$ws3.Range("F4").Value = 693
$ws3.Range("F5").Value = 279

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types - combined view)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9546
$ws4.Range("F3").Value = 2287
$ws4.Range("F4").Value = 693
$ws4.Range("F5").Value = 1272
$ws4.Range("F9").Value = 370
$ws4.Range("F10").Value = 6479
$ws4.Range("F11").Value = 279
$ws4.Range("F12").Value = 1812
$ws4.Range("F13").Value = 6382
$ws4.Range("F14").Value = 144
$ws4.Range("F15").Value = 1947
$ws4.Range("F19").Value = 37
$ws4.Range("F24").Value = 8004
$ws4.Range("F28").Value = 1753
$ws4.Range("F32").Value = 1812
$ws4.Range("F33").Value = 815
$ws4.Range("F39").Value = 16
$ws4.Range("F40").Value = 91
$ws4.Range("F45").Value = 34
